# Replace the repository URL in the final hyperlink and drop the
# hyperlink formatting (unwrap the w:hyperlink), relocating the
# trailing "_GoBack" bookmark so it sits right after the (now plain)
# URL run instead of in its own paragraph.

$d = $word.ActiveDocument

$h = $d.Hyperlinks.Item(1)

# Park the existing "_GoBack" bookmark at the end of the hyperlink's
# text *before* touching the hyperlink itself - doing this afterwards
# makes the COM host snap the bookmark back to the top of the document.
$bmPos = $h.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$bmRange.Bookmarks.Add("_GoBack")

# Unwrap the hyperlink: removes the w:hyperlink wrapper but keeps the
# inner run (and its rStyle "a8" character formatting) intact.
$h.Delete()

# Update the URL text itself (adds the "kpfu-" prefix).
$d.Content.Find.Execute(
    "https://github.com/ironsast/probability-theory-and-mathematical-statistics",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://github.com/ironsast/kpfu-probability-theory-and-mathematical-statistics",
    2) | Out-Null
